$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B18").Select()
